# Regenerate the "K" column (strikeouts, column G) of the save_data sheet.
# The sheet was previously populated from a different source ("Strike#");
# this regenerates it from the "K" source, which in general yields lower
# per-game strikeout counts than the old field. Row 1 is the header, so
# data starts at row 2 (game index 0) through row 68 (game index 66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row -> new K value (column G), taken from the
# regenerated save data.
$kValues = [ordered]@{
    2  = 2
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 3
    10 = 0
    11 = 0
    12 = 2
    13 = 0
    14 = 2
    15 = 1
    16 = 1
    17 = 0
    18 = 3
    19 = 2
    20 = 1
    21 = 3
    22 = 0
    23 = 3
    24 = 4
    25 = 2
    26 = 0
    27 = 0
    28 = 3
    29 = 1
    30 = 2
    31 = 0
    32 = 2
    33 = 0
    34 = 0
    35 = 8
    36 = 1
    37 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 1
    42 = 0
    43 = 2
    44 = 3
    45 = 0
    46 = 3
    47 = 1
    48 = 1
    49 = 5
    50 = 1
    51 = 2
    52 = 0
    53 = 0
    54 = 2
    55 = 0
    56 = 0
    57 = 1
    58 = 0
    59 = 1
    60 = 2
    61 = 1
    62 = 1
    63 = 0
    64 = 5
    65 = 5
    66 = 2
    67 = 0
    68 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
